{"js": "// Replace the 100 arithmetic-expression cells in the table with their\n// updated expressions. Each old expression is unique in the document, so we\n// locate it with a whole-word (exact) search and replace its text in place.\n// matchWholeWord avoids the one substring collision in this data set\n// (\"0+11=\" is contained inside \"10+11=\", but not as a whole-word match).\nconst replacements = [[\"9+13=\", \"99-57=\"], [\"41+20=\", \"17+81=\"], [\"45+9=\", \"53-49=\"], [\"14+79=\", \"61-10=\"], [\"30+23=\", \"54+33=\"], [\"38-37=\", \"45-19=\"], [\"88-16=\", \"62+8=\"], [\"32-14=\", \"55-19=\"], [\"23+24=\", \"71+12=\"], [\"56+36=\", \"33+16=\"], [\"84+14=\", \"56-50=\"], [\"64+11=\", \"19+18=\"], [\"26+72=\", \"87-38=\"], [\"69-15=\", \"91+6=\"], [\"80-78=\", \"89-84=\"], [\"45+49=\", \"31+56=\"], [\"8+58=\", \"84-2=\"], [\"0+11=\", \"58-29=\"], [\"28-0=\", \"31+61=\"], [\"10+11=\", \"37+0=\"], [\"95-86=\", \"96-53=\"], [\"53-7=\", \"34+12=\"], [\"47-28=\", \"62-27=\"], [\"98-10=\", \"61+6=\"], [\"63-29=\", \"78-75=\"], [\"46+6=\", \"30+68=\"], [\"52-4=\", \"41-30=\"], [\"23+49=\", \"85-10=\"], [\"79-10=\", \"25+18=\"], [\"97-1=\", \"30+19=\"], [\"2+3=\", \"56+18=\"], [\"93-41=\", \"32-19=\"], [\"17+6=\", \"33-30=\"], [\"41+46=\", \"91-33=\"], [\"46-38=\", \"87-84=\"], [\"45+34=\", \"43+30=\"], [\"85-41=\", \"10+46=\"], [\"26+61=\", \"11+36=\"], [\"55-25=\", \"89-3=\"], [\"91-75=\", \"46+20=\"], [\"12+70=\", \"54+35=\"], [\"14+81=\", \"23-1=\"], [\"94-33=\", \"66+17=\"], [\"51-44=\", \"55-41=\"], [\"44+39=\", \"68-15=\"], [\"68+16=\", \"49+34=\"], [\"57-4=\", \"54-53=\"], [\"11+49=\", \"75-14=\"], [\"57-36=\", \"43+51=\"], [\"24+47=\", \"56-19=\"], [\"88+6=\", \"43-35=\"], [\"40-16=\", \"70-22=\"], [\"5+60=\", \"6+40=\"], [\"3+14=\", \"71-22=\"], [\"44+28=\", \"80-62=\"], [\"13+61=\", \"94-24=\"], [\"8+2=\", \"84-56=\"], [\"57+0=\", \"88-77=\"], [\"71-39=\", \"33+6=\"], [\"64-62=\", \"77-18=\"], [\"56+0=\", \"44-8=\"], [\"56+29=\", \"38+45=\"], [\"21+37=\", \"27-13=\"], [\"37-25=\", \"31+60=\"], [\"95-32=\", \"68-49=\"], [\"70-66=\", \"28+52=\"], [\"52-13=\", \"46+13=\"], [\"45+39=\", \"73+14=\"], [\"82-62=\", \"46-25=\"], [\"62-23=\", \"64-40=\"], [\"21+31=\", \"83-34=\"], [\"91-48=\", \"93-59=\"], [\"29+56=\", \"90-48=\"], [\"87-6=\", \"21+43=\"], [\"45-8=\", \"29+47=\"], [\"21+48=\", \"1+96=\"], [\"61+11=\", \"84-70=\"], [\"47+46=\", \"95-75=\"], [\"77-72=\", \"82-58=\"], [\"9+1=\", \"59-38=\"], [\"66+1=\", \"56+39=\"], [\"2+10=\", \"73-64=\"], [\"7+62=\", \"17+77=\"], [\"25+41=\", \"70-14=\"], [\"49+40=\", \"57-46=\"], [\"16+78=\", \"2+79=\"], [\"85-16=\", \"83+7=\"], [\"56+20=\", \"42+53=\"], [\"60-55=\", \"69-49=\"], [\"13+16=\", \"99-54=\"], [\"65-37=\", \"74-30=\"], [\"16+18=\", \"60+36=\"], [\"83-22=\", \"1+92=\"], [\"95-65=\", \"12+71=\"], [\"48-13=\", \"1+60=\"], [\"93-88=\", \"35+7=\"], [\"80-10=\", \"95-87=\"], [\"24+58=\", \"18+11=\"], [\"64-8=\", \"26+9=\"], [\"90-40=\", \"37+17=\"]];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: true,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  // Replace every match (expected to be exactly one per string, since all\n  // 100 old values are unique in the document).\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 100 arithmetic-expression cells in the table with their\n# updated expressions. Each old expression is unique in the document, so\n# we locate/replace it with Find/Replace using MatchWholeWord to avoid the\n# one substring collision in this data set (\"0+11=\" is contained inside\n# \"10+11=\", but not as a whole-word match).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('9+13=', '99-57='),\n    @('41+20=', '17+81='),\n    @('45+9=', '53-49='),\n    @('14+79=', '61-10='),\n    @('30+23=', '54+33='),\n    @('38-37=', '45-19='),\n    @('88-16=', '62+8='),\n    @('32-14=', '55-19='),\n    @('23+24=', '71+12='),\n    @('56+36=', '33+16='),\n    @('84+14=', '56-50='),\n    @('64+11=', '19+18='),\n    @('26+72=', '87-38='),\n    @('69-15=', '91+6='),\n    @('80-78=', '89-84='),\n    @('45+49=', '31+56='),\n    @('8+58=', '84-2='),\n    @('0+11=', '58-29='),\n    @('28-0=', '31+61='),\n    @('10+11=', '37+0='),\n    @('95-86=', '96-53='),\n    @('53-7=', '34+12='),\n    @('47-28=', '62-27='),\n    @('98-10=', '61+6='),\n    @('63-29=', '78-75='),\n    @('46+6=', '30+68='),\n    @('52-4=', '41-30='),\n    @('23+49=', '85-10='),\n    @('79-10=', '25+18='),\n    @('97-1=', '30+19='),\n    @('2+3=', '56+18='),\n    @('93-41=', '32-19='),\n    @('17+6=', '33-30='),\n    @('41+46=', '91-33='),\n    @('46-38=', '87-84='),\n    @('45+34=', '43+30='),\n    @('85-41=', '10+46='),\n    @('26+61=', '11+36='),\n    @('55-25=', '89-3='),\n    @('91-75=', '46+20='),\n    @('12+70=', '54+35='),\n    @('14+81=', '23-1='),\n    @('94-33=', '66+17='),\n    @('51-44=', '55-41='),\n    @('44+39=', '68-15='),\n    @('68+16=', '49+34='),\n    @('57-4=', '54-53='),\n    @('11+49=', '75-14='),\n    @('57-36=', '43+51='),\n    @('24+47=', '56-19='),\n    @('88+6=', '43-35='),\n    @('40-16=', '70-22='),\n    @('5+60=', '6+40='),\n    @('3+14=', '71-22='),\n    @('44+28=', '80-62='),\n    @('13+61=', '94-24='),\n    @('8+2=', '84-56='),\n    @('57+0=', '88-77='),\n    @('71-39=', '33+6='),\n    @('64-62=', '77-18='),\n    @('56+0=', '44-8='),\n    @('56+29=', '38+45='),\n    @('21+37=', '27-13='),\n    @('37-25=', '31+60='),\n    @('95-32=', '68-49='),\n    @('70-66=', '28+52='),\n    @('52-13=', '46+13='),\n    @('45+39=', '73+14='),\n    @('82-62=', '46-25='),\n    @('62-23=', '64-40='),\n    @('21+31=', '83-34='),\n    @('91-48=', '93-59='),\n    @('29+56=', '90-48='),\n    @('87-6=', '21+43='),\n    @('45-8=', '29+47='),\n    @('21+48=', '1+96='),\n    @('61+11=', '84-70='),\n    @('47+46=', '95-75='),\n    @('77-72=', '82-58='),\n    @('9+1=', '59-38='),\n    @('66+1=', '56+39='),\n    @('2+10=', '73-64='),\n    @('7+62=', '17+77='),\n    @('25+41=', '70-14='),\n    @('49+40=', '57-46='),\n    @('16+78=', '2+79='),\n    @('85-16=', '83+7='),\n    @('56+20=', '42+53='),\n    @('60-55=', '69-49='),\n    @('13+16=', '99-54='),\n    @('65-37=', '74-30='),\n    @('16+18=', '60+36='),\n    @('83-22=', '1+92='),\n    @('95-65=', '12+71='),\n    @('48-13=', '1+60='),\n    @('93-88=', '35+7='),\n    @('80-10=', '95-87='),\n    @('24+58=', '18+11='),\n    @('64-8=', '26+9='),\n    @('90-40=', '37+17='),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $matched = $find.Execute(\n        [ref]$oldText,   # FindText\n        [ref]$true,      # MatchCase\n        [ref]$true,      # MatchWholeWord\n        [ref]$false,     # MatchWildcards\n        [ref]$false,     # MatchSoundsLike\n        [ref]$false,     # MatchAllWordForms\n        [ref]$true,      # Forward\n        [ref]1,          # Wrap (wdFindContinue)\n        [ref]$true,      # Format\n        [ref]$newText,   # ReplaceWith\n        [ref]2           # Replace (wdReplaceAll)\n    )\n\n    if (-not $matched) {\n        throw \"No match found for: $oldText\"\n    }\n}\n"}
